$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1359.257673344085
$ws.Range("C2").Value = 35.25400718716523
$ws.Range("D2").Value = 1107.433707964035
